$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solicitud")
$r = $ws.Range("B25:D25")
$r.Merge()
Write-Output "merged"
$r.Font.Bold = $true
Write-Output "bold set"
$r.Interior.ColorIndex = 3
Write-Output "interior set"
$r.WrapText = $true
Write-Output "wrap set"
$r.Borders.Item(7).LineStyle = 1
Write-Output "border set"
